$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel keeps them as text
$numericPriceCells = @("D5","D10","D11","D15","D20","D21","D22","D23","D26","D28","D29","D31","D33","D36","D39","D40","D43","D44","D47","D51")
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.789.97"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "1.636.29"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "215.39"
$ws.Range("E5").Value = "  -0.12%  "

$ws.Range("E6").Value = "  -0.84%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("E8").Value = "  -0.61%  "

$ws.Range("E9").Value = "  -1.18%  "

$ws.Range("D10").Value = "19.76"
$ws.Range("E10").Value = "  -2.63%  "

$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +1.19%  "

$ws.Range("E12").Value = "  +0.46%  "

$ws.Range("D13").Value = "1.863.93"
$ws.Range("E13").Value = "  -0.12%  "

$ws.Range("D14").Value = "1.637.79"
$ws.Range("E14").Value = "  -1.63%  "

$ws.Range("D15").Value = "0.562"
$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("D18").Value = "25.818.09"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").Value = "4.46"
$ws.Range("E20").Value = "  +1.23%  "

$ws.Range("D21").Value = "193.17"
$ws.Range("E21").Value = "  -0.89%  "

$ws.Range("D22").Value = "9.97"

$ws.Range("D23").Value = "6.38"
$ws.Range("E23").Value = "  +2.45%  "

$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("E25").Value = "  +2.87%  "

$ws.Range("D26").Value = "142.45"
$ws.Range("E26").Value = "  +3.04%  "

$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").Value = "6.95"
$ws.Range("E28").Value = "  +1.18%  "

$ws.Range("D29").Value = "15.54"
$ws.Range("E29").Value = "  -0.52%  "

$ws.Range("E30").Value = "  -0.60%  "

$ws.Range("D31").Value = "0.0494"
$ws.Range("E31").Value = "  -1.67%  "

$ws.Range("E32").Value = "  +0.63%  "

$ws.Range("D33").Value = "3.24"
$ws.Range("E33").Value = "  -0.49%  "

$ws.Range("E34").Value = "  -0.07%  "

$ws.Range("E35").Value = "  -0.14%  "

$ws.Range("D36").Value = "0.905"
$ws.Range("E36").Value = "  -0.41%  "

$ws.Range("D37").Value = "1.134.56"
$ws.Range("E37").Value = "  +0.37%  "

$ws.Range("E38").Value = "  -1.83%  "

$ws.Range("D39").Value = "0.545"
$ws.Range("E39").Value = "  -1.58%  "

$ws.Range("D40").Value = "0.0156"
$ws.Range("E40").Value = "  -0.97%  "

$ws.Range("E42").Value = "  +1.19%  "

$ws.Range("D43").Value = "100.51"
$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("D44").Value = "0.807"
$ws.Range("E44").Value = "  +0.60%  "

$ws.Range("D45").Value = "1.773.14"
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("D46").Value = "0.0$([char]8326)0113"
$ws.Range("E46").Value = "  +0.76%  "

$ws.Range("D47").Value = "55.31"
$ws.Range("E47").Value = "  -0.88%  "

$ws.Range("E48").Value = "  -1.55%  "

$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("E50").Value = "  +2.64%  "

$ws.Range("D51").Value = "7.51"
$ws.Range("E51").Value = "  -2.77%  "

# Clear the explicit text formatting so cell style (s attribute) matches original (default/no style)
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).ClearFormats()
}
